# Update cryptos list with latest values (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New coin name / link / price / volume(1h) values.
# A leading "'" forces values that look numeric to be stored as text,
# matching the source data which keeps these as plain strings.
$updates = @{
    'D2' = '33.583.83'
    'E2' = '  +11.04%  '
    'D3' = '1.768.82'
    'E3' = '  +5.71%  '
    'D4' = "'0.999"
    'E4' = '  +0.11%  '
    'D5' = "'229.75"
    'E5' = '  +5.28%  '
    'D6' = "'0.549"
    'E6' = '  +4.85%  '
    'E7' = '  +0.43%  '
    'D8' = "'31.45"
    'E8' = '  +6.71%  '
    'D9' = "'45.66"
    'E9' = '  +1.27%  '
    'D10' = "'0.280"
    'E10' = '  +5.28%  '
    'D11' = "'0.0665"
    'E11' = '  +7.44%  '
    'E12' = '  +2.02%  '
    'D13' = '2.026.39'
    'E13' = '  +5.78%  '
    'D14' = '1.775.36'
    'E14' = '  +6.04%  '
    'D15' = "'0.633"
    'E15' = '  +2.60%  '
    'D16' = "'10.32"
    'E16' = '  -0.85%  '
    'B17' = 'WrappedBTC'
    'C17' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    'D17' = '33.550.15'
    'E17' = '  +10.70%  '
    'B18' = 'Polkadot'
    'C18' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'D18' = "'4.31"
    'E18' = '  +7.93%  '
    'D19' = "'69.19"
    'E19' = '  +5.45%  '
    'D20' = "'260.98"
    'E20' = '  +6.30%  '
    'D21' = '0.0₃0751'
    'E21' = '  +4.94%  '
    'D22' = "'0.997"
    'E22' = '  +0.08%  '
    'E23' = '  +3.97%  '
    'D24' = "'4.38"
    'E24' = '  +1.95%  '
    'E25' = '  -0.86%  '
    'D26' = "'161.93"
    'E26' = '  +2.30%  '
    'D27' = "'16.61"
    'E27' = '  +4.38%  '
    'D28' = "'0.117"
    'E28' = '  +5.04%  '
    'D29' = "'7.11"
    'E29' = '  +5.55%  '
    'D30' = "'0.999"
    'E30' = '  +0.11%  '
    'D31' = "'3.83"
    'E31' = '  +10.61%  '
    'E32' = '  +2.99%  '
    'E33' = '  +5.91%  '
    'D34' = "'3.54"
    'E34' = '  +8.31%  '
    'D35' = '1.556.14'
    'E35' = '  +6.58%  '
    'D36' = "'1.83"
    'E36' = '  +5.32%  '
    'B37' = 'Aave'
    'C37' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D37' = "'87.56"
    'E37' = '  +9.18%  '
    'B38' = 'TrustWalletToken'
    'C38' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D38' = "'1.05"
    'E38' = '  +2.57%  '
    'D39' = "'0.629"
    'E39' = '  +7.31%  '
    'E40' = '  +4.57%  '
    'D41' = "'2.77"
    'E41' = '  +3.25%  '
    'D42' = "'2.35"
    'E42' = '  +2.94%  '
    'D43' = "'0.905"
    'E43' = '  +5.54%  '
    'D44' = "'2.12"
    'E44' = '  +6.41%  '
    'D45' = "'0.0518"
    'E45' = '  +2.94%  '
    'D46' = "'1.05"
    'E46' = '  +2.93%  '
    'B47' = 'RocketPoolETH'
    'C47' = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
    'D47' = '1.938.62'
    'E47' = '  +6.88%  '
    'B48' = 'BitcoinSV'
    'C48' = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
    'D48' = "'53.68"
    'E48' = '  +1.40%  '
    'B49' = 'FraxShare'
    'C49' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D49' = "'5.75"
    'E49' = '  +6.03%  '
    'B50' = 'PaxDollar'
    'C50' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    'D50' = "'0.999"
    'E50' = '  +0.42%  '
    'B51' = 'BabyDogeCoin'
    'C51' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D51' = '0.0₆0125'
    'E51' = '  +13.96%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
